$wb = $excel.ActiveWorkbook

# ===== Sheet: Triple Dribble =====
$ws = $wb.Worksheets.Item('Triple Dribble')
$ws.Range("A66:N66").Copy() | Out-Null
$ws.Range("A67:N67").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(67,1).Value = 'BUSTER'
$ws.Cells.Item(67,2).Value = 'WILLOW'
$ws.Cells.Item(67,3).Value = 'BEA'
$ws.Cells.Item(67,4).Value = 'BERRY'
$ws.Cells.Item(67,5).Value = 'BARLEY'
$ws.Cells.Item(67,6).Value = 'CROW'
$ws.Cells.Item(67,7).Value = 'Equipo 1'
$ws.Range("G67").Interior.Color = 16770508
$ws.Cells.Item(67,8).Value = 'TH|LeNain'
$ws.Cells.Item(67,9).Value = 'TH|iKaoss'
$ws.Cells.Item(67,10).Value = 'TH|Zhar'
$ws.Cells.Item(67,11).Value = 'NXT|amos'
$ws.Cells.Item(67,12).Value = 'NXT|Rup'
$ws.Cells.Item(67,13).Value = 'NXT|Arthur'
$ws.Cells.Item(67,14).Value = '20250725T170154.000Z'
$ws.Range("A67:N67").Copy() | Out-Null
$ws.Range("A68:N68").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(68,1).Value = 'KAZE'
$ws.Cells.Item(68,2).Value = 'CROW'
$ws.Cells.Item(68,3).Value = 'LARRY & LAWRIE'
$ws.Cells.Item(68,4).Value = 'MORTIS'
$ws.Cells.Item(68,5).Value = 'KENJI'
$ws.Cells.Item(68,6).Value = 'MEG'
$ws.Cells.Item(68,7).Value = 'Equipo 1'
$ws.Range("G68").Interior.Color = 16770508
$ws.Cells.Item(68,8).Value = 'TH|LeNain'
$ws.Cells.Item(68,9).Value = 'TH|Zhar'
$ws.Cells.Item(68,10).Value = 'TH|iKaoss'
$ws.Cells.Item(68,11).Value = 'NXT|Arthur'
$ws.Cells.Item(68,12).Value = 'NXT|amos'
$ws.Cells.Item(68,13).Value = 'NXT|Rup'
$ws.Cells.Item(68,14).Value = '20250725T165519.000Z'
$ws.Range("A68:N68").Copy() | Out-Null
$ws.Range("A69:N69").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(69,1).Value = 'KAZE'
$ws.Cells.Item(69,2).Value = 'CROW'
$ws.Cells.Item(69,3).Value = 'LARRY & LAWRIE'
$ws.Cells.Item(69,4).Value = 'MORTIS'
$ws.Cells.Item(69,5).Value = 'KENJI'
$ws.Cells.Item(69,6).Value = 'MEG'
$ws.Cells.Item(69,7).Value = 'Equipo 1'
$ws.Range("G69").Interior.Color = 16770508
$ws.Cells.Item(69,8).Value = 'TH|LeNain'
$ws.Cells.Item(69,9).Value = 'TH|Zhar'
$ws.Cells.Item(69,10).Value = 'TH|iKaoss'
$ws.Cells.Item(69,11).Value = 'NXT|Arthur'
$ws.Cells.Item(69,12).Value = 'NXT|amos'
$ws.Cells.Item(69,13).Value = 'NXT|Rup'
$ws.Cells.Item(69,14).Value = '20250725T165324.000Z'

# ===== Sheet: Sneaky Fields =====
$ws = $wb.Worksheets.Item('Sneaky Fields')
$srcWs = $wb.Worksheets.Item('Triple Dribble')
$srcWs.Range("A58:N58").Copy() | Out-Null
$ws.Range("A12:N12").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(12,1).Value = 'DRACO'
$ws.Cells.Item(12,2).Value = 'BEA'
$ws.Cells.Item(12,3).Value = 'CORDELIUS'
$ws.Cells.Item(12,4).Value = 'ALLI'
$ws.Cells.Item(12,5).Value = 'CHESTER'
$ws.Cells.Item(12,6).Value = 'MEEPLE'
$ws.Cells.Item(12,7).Value = 'Equipo 1'
$ws.Range("G12").Interior.Color = 16770508
$ws.Cells.Item(12,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(12,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(12,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(12,11).Value = 'Bielz'
$ws.Cells.Item(12,12).Value = 'GO|Yichy❦'
$ws.Cells.Item(12,13).Value = 'Tilo🍥'
$ws.Cells.Item(12,14).Value = '20250725T165057.000Z'
$ws.Range("A12:N12").Copy() | Out-Null
$ws.Range("A13:N13").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(13,1).Value = 'DRACO'
$ws.Cells.Item(13,2).Value = 'BEA'
$ws.Cells.Item(13,3).Value = 'CORDELIUS'
$ws.Cells.Item(13,4).Value = 'ALLI'
$ws.Cells.Item(13,5).Value = 'CHESTER'
$ws.Cells.Item(13,6).Value = 'MEEPLE'
$ws.Cells.Item(13,7).Value = 'Equipo 1'
$ws.Range("G13").Interior.Color = 16770508
$ws.Cells.Item(13,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(13,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(13,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(13,11).Value = 'Bielz'
$ws.Cells.Item(13,12).Value = 'GO|Yichy❦'
$ws.Cells.Item(13,13).Value = 'Tilo🍥'
$ws.Cells.Item(13,14).Value = '20250725T164844.000Z'
$ws.Range("A13:N13").Copy() | Out-Null
$ws.Range("A14:N14").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14,1).Value = 'CORDELIUS'
$ws.Cells.Item(14,2).Value = 'SPIKE'
$ws.Cells.Item(14,3).Value = 'MEG'
$ws.Cells.Item(14,4).Value = 'R-T'
$ws.Cells.Item(14,5).Value = 'FINX'
$ws.Cells.Item(14,6).Value = 'MOE'
$ws.Cells.Item(14,7).Value = 'Equipo 1'
$ws.Range("G14").Interior.Color = 16770508
$ws.Cells.Item(14,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(14,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(14,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(14,11).Value = 'Bielz'
$ws.Cells.Item(14,12).Value = 'Tilo🍥'
$ws.Cells.Item(14,13).Value = 'GO|Yichy❦'
$ws.Cells.Item(14,14).Value = '20250725T164220.000Z'
$ws.Range("A14:N14").Copy() | Out-Null
$ws.Range("A15:N15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,1).Value = 'CORDELIUS'
$ws.Cells.Item(15,2).Value = 'SPIKE'
$ws.Cells.Item(15,3).Value = 'MEG'
$ws.Cells.Item(15,4).Value = 'R-T'
$ws.Cells.Item(15,5).Value = 'FINX'
$ws.Cells.Item(15,6).Value = 'MOE'
$ws.Cells.Item(15,7).Value = 'Equipo 1'
$ws.Range("G15").Interior.Color = 16770508
$ws.Cells.Item(15,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(15,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(15,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(15,11).Value = 'Bielz'
$ws.Cells.Item(15,12).Value = 'Tilo🍥'
$ws.Cells.Item(15,13).Value = 'GO|Yichy❦'
$ws.Cells.Item(15,14).Value = '20250725T164014.000Z'
$ws.Range("A15:N15").Copy() | Out-Null
$ws.Range("A16:N16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16,1).Value = 'CORDELIUS'
$ws.Cells.Item(16,2).Value = 'SPIKE'
$ws.Cells.Item(16,3).Value = 'MEG'
$ws.Cells.Item(16,4).Value = 'R-T'
$ws.Cells.Item(16,5).Value = 'FINX'
$ws.Cells.Item(16,6).Value = 'MOE'
$ws.Cells.Item(16,7).Value = 'Empate'
$ws.Range("G16").Interior.Color = 14277081
$ws.Cells.Item(16,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(16,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(16,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(16,11).Value = 'Bielz'
$ws.Cells.Item(16,12).Value = 'Tilo🍥'
$ws.Cells.Item(16,13).Value = 'GO|Yichy❦'
$ws.Cells.Item(16,14).Value = '20250725T163753.000Z'
$ws.Range("A16:N16").Copy() | Out-Null
$ws.Range("A17:N17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17,1).Value = 'CORDELIUS'
$ws.Cells.Item(17,2).Value = 'SPIKE'
$ws.Cells.Item(17,3).Value = 'MEG'
$ws.Cells.Item(17,4).Value = 'R-T'
$ws.Cells.Item(17,5).Value = 'FINX'
$ws.Cells.Item(17,6).Value = 'MOE'
$ws.Cells.Item(17,7).Value = 'Equipo 2'
$ws.Range("G17").Interior.Color = 13421812
$ws.Cells.Item(17,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(17,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(17,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(17,11).Value = 'Bielz'
$ws.Cells.Item(17,12).Value = 'Tilo🍥'
$ws.Cells.Item(17,13).Value = 'GO|Yichy❦'
$ws.Cells.Item(17,14).Value = '20250725T163402.000Z'

# ===== Sheet: Hot Potato =====
$ws = $wb.Worksheets.Item('Hot Potato')
$ws.Range("A69:N69").Copy() | Out-Null
$ws.Range("A70:N70").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(70,1).Value = 'CROW'
$ws.Cells.Item(70,2).Value = 'CHUCK'
$ws.Cells.Item(70,3).Value = 'GRIFF'
$ws.Cells.Item(70,4).Value = 'KAZE'
$ws.Cells.Item(70,5).Value = 'AMBER'
$ws.Cells.Item(70,6).Value = 'CHARLIE'
$ws.Cells.Item(70,7).Value = 'Equipo 2'
$ws.Range("G70").Interior.Color = 13421812
$ws.Cells.Item(70,8).Value = 'IC|Mebius'
$ws.Cells.Item(70,9).Value = 'IC|RamaZR'
$ws.Cells.Item(70,10).Value = 'IC|Nob?'
$ws.Cells.Item(70,11).Value = 'FUT|Nowy297'
$ws.Cells.Item(70,12).Value = 'FUT|MeOw'
$ws.Cells.Item(70,13).Value = 'FUT|GeRo'
$ws.Cells.Item(70,14).Value = '20250725T170221.000Z'
$ws.Range("A70:N70").Copy() | Out-Null
$ws.Range("A71:N71").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(71,1).Value = 'LILY'
$ws.Cells.Item(71,2).Value = 'MICO'
$ws.Cells.Item(71,3).Value = 'LOU'
$ws.Cells.Item(71,4).Value = 'LUMI'
$ws.Cells.Item(71,5).Value = 'BULL'
$ws.Cells.Item(71,6).Value = 'SHADE'
$ws.Cells.Item(71,7).Value = 'Equipo 2'
$ws.Range("G71").Interior.Color = 13421812
$ws.Cells.Item(71,8).Value = 'IC|Mebius'
$ws.Cells.Item(71,9).Value = 'IC|RamaZR'
$ws.Cells.Item(71,10).Value = 'IC|Nob?'
$ws.Cells.Item(71,11).Value = 'FUT|GeRo'
$ws.Cells.Item(71,12).Value = 'FUT|Nowy297'
$ws.Cells.Item(71,13).Value = 'FUT|MeOw'
$ws.Cells.Item(71,14).Value = '20250725T165545.000Z'
$ws.Range("A71:N71").Copy() | Out-Null
$ws.Range("A72:N72").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(72,1).Value = 'LILY'
$ws.Cells.Item(72,2).Value = 'MICO'
$ws.Cells.Item(72,3).Value = 'LOU'
$ws.Cells.Item(72,4).Value = 'LUMI'
$ws.Cells.Item(72,5).Value = 'BULL'
$ws.Cells.Item(72,6).Value = 'SHADE'
$ws.Cells.Item(72,7).Value = 'Equipo 2'
$ws.Range("G72").Interior.Color = 13421812
$ws.Cells.Item(72,8).Value = 'IC|Mebius'
$ws.Cells.Item(72,9).Value = 'IC|RamaZR'
$ws.Cells.Item(72,10).Value = 'IC|Nob?'
$ws.Cells.Item(72,11).Value = 'FUT|GeRo'
$ws.Cells.Item(72,12).Value = 'FUT|Nowy297'
$ws.Cells.Item(72,13).Value = 'FUT|MeOw'
$ws.Cells.Item(72,14).Value = '20250725T165430.000Z'
$ws.Range("A72:N72").Copy() | Out-Null
$ws.Range("A73:N73").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(73,1).Value = 'LILY'
$ws.Cells.Item(73,2).Value = 'MICO'
$ws.Cells.Item(73,3).Value = 'LOU'
$ws.Cells.Item(73,4).Value = 'LUMI'
$ws.Cells.Item(73,5).Value = 'BULL'
$ws.Cells.Item(73,6).Value = 'SHADE'
$ws.Cells.Item(73,7).Value = 'Equipo 1'
$ws.Range("G73").Interior.Color = 16770508
$ws.Cells.Item(73,8).Value = 'IC|Mebius'
$ws.Cells.Item(73,9).Value = 'IC|RamaZR'
$ws.Cells.Item(73,10).Value = 'IC|Nob?'
$ws.Cells.Item(73,11).Value = 'FUT|GeRo'
$ws.Cells.Item(73,12).Value = 'FUT|Nowy297'
$ws.Cells.Item(73,13).Value = 'FUT|MeOw'
$ws.Cells.Item(73,14).Value = '20250725T165301.000Z'
$ws.Range("A73:N73").Copy() | Out-Null
$ws.Range("A74:N74").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(74,1).Value = 'CROW'
$ws.Cells.Item(74,2).Value = 'CHUCK'
$ws.Cells.Item(74,3).Value = 'GRIFF'
$ws.Cells.Item(74,4).Value = 'KAZE'
$ws.Cells.Item(74,5).Value = 'AMBER'
$ws.Cells.Item(74,6).Value = 'CHARLIE'
$ws.Cells.Item(74,7).Value = 'Equipo 1'
$ws.Range("G74").Interior.Color = 16770508
$ws.Cells.Item(74,8).Value = 'IC|Mebius'
$ws.Cells.Item(74,9).Value = 'IC|RamaZR'
$ws.Cells.Item(74,10).Value = 'IC|Nob?'
$ws.Cells.Item(74,11).Value = 'FUT|Nowy297'
$ws.Cells.Item(74,12).Value = 'FUT|MeOw'
$ws.Cells.Item(74,13).Value = 'FUT|GeRo'
$ws.Cells.Item(74,14).Value = '20250725T170407.000Z'

# ===== Sheet: Layer Cake =====
$ws = $wb.Worksheets.Item('Layer Cake')
$ws.Range("A75:N75").Copy() | Out-Null
$ws.Range("A76:N76").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(76,1).Value = 'DOUG'
$ws.Cells.Item(76,2).Value = 'GENE'
$ws.Cells.Item(76,3).Value = 'JANET'
$ws.Cells.Item(76,4).Value = 'GUS'
$ws.Cells.Item(76,5).Value = 'HANK'
$ws.Cells.Item(76,6).Value = 'KIT'
$ws.Cells.Item(76,7).Value = 'Equipo 2'
$ws.Range("G76").Interior.Color = 13421812
$ws.Cells.Item(76,8).Value = 'IC|RamaZR'
$ws.Cells.Item(76,9).Value = 'IC|Nob?'
$ws.Cells.Item(76,10).Value = 'IC|Mebius'
$ws.Cells.Item(76,11).Value = 'FUT|GeRo'
$ws.Cells.Item(76,12).Value = 'FUT|Nowy297'
$ws.Cells.Item(76,13).Value = 'FUT|MeOw'
$ws.Cells.Item(76,14).Value = '20250725T164827.000Z'
$ws.Range("A76:N76").Copy() | Out-Null
$ws.Range("A77:N77").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(77,1).Value = 'DOUG'
$ws.Cells.Item(77,2).Value = 'GENE'
$ws.Cells.Item(77,3).Value = 'JANET'
$ws.Cells.Item(77,4).Value = 'GUS'
$ws.Cells.Item(77,5).Value = 'HANK'
$ws.Cells.Item(77,6).Value = 'KIT'
$ws.Cells.Item(77,7).Value = 'Equipo 2'
$ws.Range("G77").Interior.Color = 13421812
$ws.Cells.Item(77,8).Value = 'IC|RamaZR'
$ws.Cells.Item(77,9).Value = 'IC|Nob?'
$ws.Cells.Item(77,10).Value = 'IC|Mebius'
$ws.Cells.Item(77,11).Value = 'FUT|GeRo'
$ws.Cells.Item(77,12).Value = 'FUT|Nowy297'
$ws.Cells.Item(77,13).Value = 'FUT|MeOw'
$ws.Cells.Item(77,14).Value = '20250725T164620.000Z'
$ws.Range("A77:N77").Copy() | Out-Null
$ws.Range("A78:N78").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(78,1).Value = 'CORDELIUS'
$ws.Cells.Item(78,2).Value = 'SQUEAK'
$ws.Cells.Item(78,3).Value = 'GUS'
$ws.Cells.Item(78,4).Value = 'CROW'
$ws.Cells.Item(78,5).Value = 'DOUG'
$ws.Cells.Item(78,6).Value = 'JANET'
$ws.Cells.Item(78,7).Value = 'Equipo 1'
$ws.Range("G78").Interior.Color = 16770508
$ws.Cells.Item(78,8).Value = 'IC|Mebius'
$ws.Cells.Item(78,9).Value = 'IC|RamaZR'
$ws.Cells.Item(78,10).Value = 'IC|Nob?'
$ws.Cells.Item(78,11).Value = 'FUT|GeRo'
$ws.Cells.Item(78,12).Value = 'FUT|Nowy297'
$ws.Cells.Item(78,13).Value = 'FUT|MeOw'
$ws.Cells.Item(78,14).Value = '20250725T163953.000Z'
$ws.Range("A78:N78").Copy() | Out-Null
$ws.Range("A79:N79").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(79,1).Value = 'CORDELIUS'
$ws.Cells.Item(79,2).Value = 'SQUEAK'
$ws.Cells.Item(79,3).Value = 'GUS'
$ws.Cells.Item(79,4).Value = 'CROW'
$ws.Cells.Item(79,5).Value = 'DOUG'
$ws.Cells.Item(79,6).Value = 'JANET'
$ws.Cells.Item(79,7).Value = 'Equipo 1'
$ws.Range("G79").Interior.Color = 16770508
$ws.Cells.Item(79,8).Value = 'IC|Mebius'
$ws.Cells.Item(79,9).Value = 'IC|RamaZR'
$ws.Cells.Item(79,10).Value = 'IC|Nob?'
$ws.Cells.Item(79,11).Value = 'FUT|GeRo'
$ws.Cells.Item(79,12).Value = 'FUT|Nowy297'
$ws.Cells.Item(79,13).Value = 'FUT|MeOw'
$ws.Cells.Item(79,14).Value = '20250725T163806.000Z'

# ===== Sheet: Dry Season =====
$ws = $wb.Worksheets.Item('Dry Season')
$ws.Range("A44:N44").Copy() | Out-Null
$ws.Range("A45:N45").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(45,1).Value = 'CARL'
$ws.Cells.Item(45,2).Value = 'GUS'
$ws.Cells.Item(45,3).Value = 'ANGELO'
$ws.Cells.Item(45,4).Value = 'R-T'
$ws.Cells.Item(45,5).Value = 'MAX'
$ws.Cells.Item(45,6).Value = 'MR. P'
$ws.Cells.Item(45,7).Value = 'Equipo 2'
$ws.Range("G45").Interior.Color = 13421812
$ws.Cells.Item(45,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(45,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(45,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(45,11).Value = 'Bielz'
$ws.Cells.Item(45,12).Value = 'Tilo🍥'
$ws.Cells.Item(45,13).Value = 'GO|Yichy❦'
$ws.Cells.Item(45,14).Value = '20250725T170152.000Z'
$ws.Range("A45:N45").Copy() | Out-Null
$ws.Range("A46:N46").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(46,1).Value = 'CARL'
$ws.Cells.Item(46,2).Value = 'GUS'
$ws.Cells.Item(46,3).Value = 'ANGELO'
$ws.Cells.Item(46,4).Value = 'R-T'
$ws.Cells.Item(46,5).Value = 'MAX'
$ws.Cells.Item(46,6).Value = 'MR. P'
$ws.Cells.Item(46,7).Value = 'Equipo 2'
$ws.Range("G46").Interior.Color = 13421812
$ws.Cells.Item(46,8).Value = 'LOUD|FireCrow'
$ws.Cells.Item(46,9).Value = 'LOUD|Edinho'
$ws.Cells.Item(46,10).Value = 'LOUD|KaioDog'
$ws.Cells.Item(46,11).Value = 'Bielz'
$ws.Cells.Item(46,12).Value = 'Tilo🍥'
$ws.Cells.Item(46,13).Value = 'GO|Yichy❦'
$ws.Cells.Item(46,14).Value = '20250725T165931.000Z'

# ===== Sheet: Pit Stop =====
$ws = $wb.Worksheets.Item('Pit Stop')
$ws.Range("A56:N56").Copy() | Out-Null
$ws.Range("A57:N57").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(57,1).Value = 'LOU'
$ws.Cells.Item(57,2).Value = 'MICO'
$ws.Cells.Item(57,3).Value = 'BERRY'
$ws.Cells.Item(57,4).Value = 'KIT'
$ws.Cells.Item(57,5).Value = 'KAZE'
$ws.Cells.Item(57,6).Value = 'BULL'
$ws.Cells.Item(57,7).Value = 'Equipo 1'
$ws.Range("G57").Interior.Color = 16770508
$ws.Cells.Item(57,8).Value = 'TH|Zhar'
$ws.Cells.Item(57,9).Value = 'TH|LeNain'
$ws.Cells.Item(57,10).Value = 'TH|iKaoss'
$ws.Cells.Item(57,11).Value = 'NXT|Rup'
$ws.Cells.Item(57,12).Value = 'NXT|Arthur'
$ws.Cells.Item(57,13).Value = 'NXT|amos'
$ws.Cells.Item(57,14).Value = '20250725T164737.000Z'
$ws.Range("A57:N57").Copy() | Out-Null
$ws.Range("A58:N58").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(58,1).Value = 'LOU'
$ws.Cells.Item(58,2).Value = 'MICO'
$ws.Cells.Item(58,3).Value = 'BERRY'
$ws.Cells.Item(58,4).Value = 'KIT'
$ws.Cells.Item(58,5).Value = 'KAZE'
$ws.Cells.Item(58,6).Value = 'BULL'
$ws.Cells.Item(58,7).Value = 'Equipo 2'
$ws.Range("G58").Interior.Color = 13421812
$ws.Cells.Item(58,8).Value = 'TH|Zhar'
$ws.Cells.Item(58,9).Value = 'TH|LeNain'
$ws.Cells.Item(58,10).Value = 'TH|iKaoss'
$ws.Cells.Item(58,11).Value = 'NXT|Rup'
$ws.Cells.Item(58,12).Value = 'NXT|Arthur'
$ws.Cells.Item(58,13).Value = 'NXT|amos'
$ws.Cells.Item(58,14).Value = '20250725T164601.000Z'
$ws.Range("A58:N58").Copy() | Out-Null
$ws.Range("A59:N59").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(59,1).Value = 'LOU'
$ws.Cells.Item(59,2).Value = 'MICO'
$ws.Cells.Item(59,3).Value = 'BERRY'
$ws.Cells.Item(59,4).Value = 'KIT'
$ws.Cells.Item(59,5).Value = 'KAZE'
$ws.Cells.Item(59,6).Value = 'BULL'
$ws.Cells.Item(59,7).Value = 'Equipo 1'
$ws.Range("G59").Interior.Color = 16770508
$ws.Cells.Item(59,8).Value = 'TH|Zhar'
$ws.Cells.Item(59,9).Value = 'TH|LeNain'
$ws.Cells.Item(59,10).Value = 'TH|iKaoss'
$ws.Cells.Item(59,11).Value = 'NXT|Rup'
$ws.Cells.Item(59,12).Value = 'NXT|Arthur'
$ws.Cells.Item(59,13).Value = 'NXT|amos'
$ws.Cells.Item(59,14).Value = '20250725T164442.000Z'
$ws.Range("A59:N59").Copy() | Out-Null
$ws.Range("A60:N60").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(60,1).Value = 'BUZZ'
$ws.Cells.Item(60,2).Value = 'MICO'
$ws.Cells.Item(60,3).Value = 'CHARLIE'
$ws.Cells.Item(60,4).Value = 'KIT'
$ws.Cells.Item(60,5).Value = 'CARL'
$ws.Cells.Item(60,6).Value = 'KAZE'
$ws.Cells.Item(60,7).Value = 'Equipo 2'
$ws.Range("G60").Interior.Color = 13421812
$ws.Cells.Item(60,8).Value = 'TH|LeNain'
$ws.Cells.Item(60,9).Value = 'TH|iKaoss'
$ws.Cells.Item(60,10).Value = 'TH|Zhar'
$ws.Cells.Item(60,11).Value = 'NXT|Rup'
$ws.Cells.Item(60,12).Value = 'NXT|amos'
$ws.Cells.Item(60,13).Value = 'NXT|Arthur'
$ws.Cells.Item(60,14).Value = '20250725T164035.000Z'
$ws.Range("A60:N60").Copy() | Out-Null
$ws.Range("A61:N61").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(61,1).Value = 'BUZZ'
$ws.Cells.Item(61,2).Value = 'MICO'
$ws.Cells.Item(61,3).Value = 'CHARLIE'
$ws.Cells.Item(61,4).Value = 'KIT'
$ws.Cells.Item(61,5).Value = 'CARL'
$ws.Cells.Item(61,6).Value = 'KAZE'
$ws.Cells.Item(61,7).Value = 'Equipo 2'
$ws.Range("G61").Interior.Color = 13421812
$ws.Cells.Item(61,8).Value = 'TH|LeNain'
$ws.Cells.Item(61,9).Value = 'TH|iKaoss'
$ws.Cells.Item(61,10).Value = 'TH|Zhar'
$ws.Cells.Item(61,11).Value = 'NXT|Rup'
$ws.Cells.Item(61,12).Value = 'NXT|amos'
$ws.Cells.Item(61,13).Value = 'NXT|Arthur'
$ws.Cells.Item(61,14).Value = '20250725T163830.000Z'
$ws.Range("A61:N61").Copy() | Out-Null
$ws.Range("A62:N62").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(62,1).Value = 'BUZZ'
$ws.Cells.Item(62,2).Value = 'MICO'
$ws.Cells.Item(62,3).Value = 'CHARLIE'
$ws.Cells.Item(62,4).Value = 'KIT'
$ws.Cells.Item(62,5).Value = 'CARL'
$ws.Cells.Item(62,6).Value = 'KAZE'
$ws.Cells.Item(62,7).Value = 'Equipo 1'
$ws.Range("G62").Interior.Color = 16770508
$ws.Cells.Item(62,8).Value = 'TH|LeNain'
$ws.Cells.Item(62,9).Value = 'TH|iKaoss'
$ws.Cells.Item(62,10).Value = 'TH|Zhar'
$ws.Cells.Item(62,11).Value = 'NXT|Rup'
$ws.Cells.Item(62,12).Value = 'NXT|amos'
$ws.Cells.Item(62,13).Value = 'NXT|Arthur'
$ws.Cells.Item(62,14).Value = '20250725T163617.000Z'

$excel.CutCopyMode = 0
Write-Host "All rows inserted."